# Apply the "2023-10-01sales" update:
#  - Rewrite the data rows (2-6) with the new sales records.
#  - Drop the old trailing rows (7-10) so the table shrinks to A1:E6.
# Values are written as TEXT (shared strings), matching the source file,
# by forcing a text number format before the write and then resetting the
# cell style back to Normal so no stray per-cell style sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: PC 121212 | 2023-10-30 | Caro White Cream | 1800 | 1
Set-TextValue $ws.Range("A2") "121212"
Set-TextValue $ws.Range("B2") "2023-10-30"
Set-TextValue $ws.Range("C2") "Caro White Cream"
Set-TextValue $ws.Range("D2") "1800"
Set-TextValue $ws.Range("E2") "1"

# Row 3: PC 121212 | 2023-10-30 | Caro White Cream | 1800 | 1
Set-TextValue $ws.Range("A3") "121212"
Set-TextValue $ws.Range("B3") "2023-10-30"
Set-TextValue $ws.Range("C3") "Caro White Cream"
Set-TextValue $ws.Range("D3") "1800"
Set-TextValue $ws.Range("E3") "1"

# Row 4: PC 121210 | 2023-10-30 | Oral B | 800 | 1
Set-TextValue $ws.Range("A4") "121210"
Set-TextValue $ws.Range("B4") "2023-10-30"
Set-TextValue $ws.Range("C4") "Oral B"
Set-TextValue $ws.Range("D4") "800"
Set-TextValue $ws.Range("E4") "1"

# Row 5: PC 121212 | 2023-10-31 | Caro White Cream | 1800 | 1
Set-TextValue $ws.Range("A5") "121212"
Set-TextValue $ws.Range("B5") "2023-10-31"
Set-TextValue $ws.Range("C5") "Caro White Cream"
Set-TextValue $ws.Range("D5") "1800"
Set-TextValue $ws.Range("E5") "1"

# Row 6: PC 121212 | 2023-10-31 | Caro White Cream | 1800 | 1
Set-TextValue $ws.Range("A6") "121212"
Set-TextValue $ws.Range("B6") "2023-10-31"
Set-TextValue $ws.Range("C6") "Caro White Cream"
Set-TextValue $ws.Range("D6") "1800"
Set-TextValue $ws.Range("E6") "1"

# Remove the now-unused trailing rows 7-10.
$ws.Rows("7:10").Delete()
